# Update the "Förändrad" (Changed) date column (C) for every data row
# from 2023-09-21 (45190) to 2023-09-23 (45192).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$xlUp = -4162

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End($xlUp).Row
if ($lastRow -lt 2) {
    $lastRow = $ws.UsedRange.Rows.Count
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45190) {
        $cell.Value2 = 45192
    }
}
